$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired content for rows 4-14 (product, process, location, database)
$data = @(
    @(4,  "transport, passenger car, electric", "transport, passenger car, electric", "GLO", "ecoinvent38_cutoff"),
    @(5,  "transport, passenger car, small size, petrol, EURO 5", "transport, passenger car, small size, petrol, EURO 5", "RER", "ecoinvent38_cutoff"),
    @(6,  "electricity, high voltage", "treatment of bagasse, from sweet sorghum, in heat and power co-generation unit, 6400kW thermal", "GLO", "ecoinvent38_cutoff"),
    @(7,  "heat, district or industrial, other than natural gas", "treatment of bagasse, from sweet sorghum, in heat and power co-generation unit, 6400kW thermal", "GLO", "ecoinvent38_cutoff"),
    @(8,  "electricity, low voltage", "wood pellets, burned in stirling heat and power co-generation unit, 3kW electrical, future", "CH", "ecoinvent38_cutoff"),
    @(9,  "act1", "prod1", "DE", "eco_export_SS"),
    @(10, "act2", "prod2", "DE", "eco_export_SS"),
    @(11, "act3", "prod3", "DE", "eco_export_SS"),
    @(12, "act1", "prod1", "DE", "eco_export_SS"),
    @(13, "act2", "prod2", "DE", "eco_export_SS"),
    @(14, "act3", "prod3", "DE", "eco_export_SS")
)

# Fill columns C and D first (these only reuse already-existing shared strings)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Then column B (process / product names) so the new "prod1..prod3" strings
# are appended to the shared string table before the "act1..act3" ones,
# matching the original authoring order.
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

# Finally column A (activity names)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
}

# Update the active selection to match the target file
$ws.Range("D17").Select()

$wb.Save()
